$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2022-05-30 T 18:36:33 UTC"
$ws.Range("B2").Value = 1905.94020055625
$ws.Range("C2").Value = 0.7886609999999999
$ws.Range("D2").Value = 1.259291

$ws.Range("A3").Value = "2022-05-30 T 18:36:33 UTC"
$ws.Range("B3").Value = 1905.94020055625
$ws.Range("C3").Value = 0.7886609999999999
$ws.Range("D3").Value = 1.259291
